$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("A9").Value = 111675587
$ws.Range("I9").Value = "'3"
$ws.Range("I9").ClearFormats()
$ws.Range("Q9").Value = 690345
$ws.Range("R9").Value = 6661441
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()

# --- Row 10 ---
$ws.Range("A10").Value = 111675585
$ws.Range("I10").Value = "'1"
$ws.Range("I10").ClearFormats()
$ws.Range("Q10").Value = 690350
$ws.Range("R10").Value = 6661440
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()

# --- Row 11 ---
$ws.Range("A11").Value = 111675586
$ws.Range("I11").Value = "'2"
$ws.Range("I11").ClearFormats()
$ws.Range("Q11").Value = 690349
$ws.Range("R11").Value = 6661441
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()
